$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$jobDescription = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"

$ws.Range("A4").Value = "JD_003"
$ws.Range("B4").Value = "Junior RPA Developer"
$ws.Range("C4").Value = $jobDescription
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
